$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11 (sheet ALC)
$ws.Range("H11").Value = 206.375
$ws.Range("I11").Value = 206.375
$ws.Range("K11").Value = 206.375
$ws.Range("M11").Value = -66.375
# Row 40 (sheet ALC)
$ws.Range("H40").Value = 1954.2858
$ws.Range("I40").Value = 1977.5
$ws.Range("J40").Value = 1945
$ws.Range("K40").Value = 1977.5
$ws.Range("L40").Value = 1945
$ws.Range("M40").Value = -1802.5
$ws.Range("N40").Value = -2295
# Row 88 (sheet ALC)
$ws.Range("H88").Value = 1398.25
$ws.Range("I88").Value = 793
$ws.Range("J88").Value = 1600
$ws.Range("K88").Value = 793
$ws.Range("L88").Value = 1600
$ws.Range("M88").Value = -387
$ws.Range("N88").Value = -2412
# Row 91 (sheet ALC)
$ws.Range("H91").Value = 1398.25
$ws.Range("I91").Value = 793
$ws.Range("J91").Value = 1600
$ws.Range("K91").Value = 793
$ws.Range("L91").Value = 1600
$ws.Range("M91").Value = 611
$ws.Range("N91").Value = -4408
# Row 107 (sheet ALC)
$ws.Range("H107").Value = 5375
$ws.Range("I107").Value = 750
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 750
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = 1170
$ws.Range("N107").Value = -13840
# Row 116 (sheet ALC)
$ws.Range("H116").Value = 3074.8
$ws.Range("I116").Value = 3074.8
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 3074.8
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 367.1999999999998
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (sheet ARM)
$ws.Range("H2").Value = 1040.6666
$ws.Range("I2").Value = 981
$ws.Range("K2").Value = 981
$ws.Range("M2").Value = -868
# Row 32 (sheet ARM)
$ws.Range("H32").Value = 2147.4583
$ws.Range("I32").Value = 2147.4583
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2147.4583
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1860.4583
$ws.Range("N32").ClearContents()
# Row 74 (sheet ARM)
$ws.Range("H74").Value = 3625.111
$ws.Range("J74").Value = 2470
$ws.Range("L74").Value = 2470
$ws.Range("N74").Value = -4218
# Row 77 (sheet ARM)
$ws.Range("H77").Value = 3625.111
$ws.Range("J77").Value = 2470
$ws.Range("L77").Value = 12350
$ws.Range("N77").Value = -21086
# Row 88 (sheet ARM)
$ws.Range("H88").Value = 3116.5
$ws.Range("I88").Value = 1110
$ws.Range("J88").Value = 3663.7273
$ws.Range("K88").Value = 1110
$ws.Range("L88").Value = 3663.7273
$ws.Range("M88").Value = -704
$ws.Range("N88").Value = -4475.7273
# Row 91 (sheet ARM)
$ws.Range("H91").Value = 3116.5
$ws.Range("I91").Value = 1110
$ws.Range("J91").Value = 3663.7273
$ws.Range("K91").Value = 1110
$ws.Range("L91").Value = 3663.7273
$ws.Range("M91").Value = 294
$ws.Range("N91").Value = -6471.7273
# Row 116 (sheet ARM)
$ws.Range("H116").Value = 1040.6666
$ws.Range("I116").Value = 981
$ws.Range("K116").Value = 981
$ws.Range("M116").Value = 1313

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (sheet BSM)
$ws.Range("H3").Value = 1040.6666
$ws.Range("I3").Value = 981
$ws.Range("K3").Value = 981
$ws.Range("M3").Value = -867
# Row 22 (sheet BSM)
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 50 (sheet CRP)
$ws.Range("H50").Value = 28180.5
$ws.Range("J50").Value = 29600
$ws.Range("L50").Value = 29600
$ws.Range("N50").Value = -30850
# Row 60 (sheet CRP)
$ws.Range("H60").Value = 22250
$ws.Range("I60").Value = 14000
$ws.Range("K60").Value = 14000
$ws.Range("M60").Value = -13489

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (sheet CUL)
$ws.Range("H4").Value = 2098.75
$ws.Range("I4").Value = 2592.8572
$ws.Range("K4").Value = 7778.571599999999
$ws.Range("M4").Value = -7666.571599999999

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (sheet GSM)
$ws.Range("H2").Value = 151.75
$ws.Range("I2").Value = 165.16667
$ws.Range("K2").Value = 165.16667
$ws.Range("M2").Value = -52.16667000000001
# Row 63 (sheet GSM)
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
# Row 66 (sheet GSM)
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()
# Row 99 (sheet GSM)
$ws.Range("H99").Value = 6820.5
$ws.Range("I99").Value = 6820.5
$ws.Range("K99").Value = 6820.5
$ws.Range("M99").Value = -4574.5
# Row 102 (sheet GSM)
$ws.Range("H102").Value = 4129.3335
$ws.Range("I102").Value = 4129.3335
$ws.Range("K102").Value = 4129.3335
$ws.Range("M102").Value = -2507.3335
# Row 107 (sheet GSM)
$ws.Range("H107").Value = 297.75
$ws.Range("I107").Value = 297.75
$ws.Range("K107").Value = 297.75
$ws.Range("M107").Value = 1622.25
# Row 113 (sheet GSM)
$ws.Range("H113").Value = 3483.6667
$ws.Range("J113").Value = 2500
$ws.Range("L113").Value = 2500
$ws.Range("N113").Value = -6840
# Row 122 (sheet GSM)
$ws.Range("H122").Value = 4145.25
$ws.Range("I122").Value = 2619
$ws.Range("J122").Value = 5061
$ws.Range("K122").Value = 7857
$ws.Range("L122").Value = 15183
$ws.Range("M122").Value = -5407
$ws.Range("N122").Value = -20083

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (sheet LTW)
$ws.Range("H16").Value = 3455.8
$ws.Range("J16").Value = 3290
$ws.Range("L16").Value = 3290
$ws.Range("N16").Value = -3630
# Row 46 (sheet LTW)
$ws.Range("H46").Value = 3461.074
$ws.Range("I46").Value = 2749.923
$ws.Range("J46").Value = 4121.4287
$ws.Range("K46").Value = 2749.923
$ws.Range("L46").Value = 4121.4287
$ws.Range("M46").Value = -2561.923
$ws.Range("N46").Value = -4497.4287
# Row 132 (sheet LTW)
$ws.Range("H132").Value = 6217.45
$ws.Range("I132").Value = 4939.125
$ws.Range("K132").Value = 14817.375
$ws.Range("M132").Value = -12287.375

$ws = $wb.Worksheets.Item("WVR")
# Row 126 (sheet WVR)
$ws.Range("H126").Value = 1384.6666
$ws.Range("I126").Value = 1402
$ws.Range("K126").Value = 4206
$ws.Range("M126").Value = -1736
# Row 136 (sheet WVR)
$ws.Range("H136").Value = 7147.2856
$ws.Range("I136").Value = 6796.636
$ws.Range("J136").Value = 8433
$ws.Range("K136").Value = 20389.908
$ws.Range("L136").Value = 25299
$ws.Range("M136").Value = -17839.908
$ws.Range("N136").Value = -30399
